$d = $word.ActiveDocument

# --- 1. Update the "GA1-220501093-AA1" heading paragraph ---
# Find the paragraph that currently reads "GA1-220501093-AA1" (centered,
# firstLine indent only) and:
#   - change its text to the new activity code
#   - drop the centered justification
#   - add a left indent of 3528 twips (176.4 pt), keeping the firstLine indent
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "GA1-220501093-AA1") {
        $target = $p
        break
    }
}

$target.Range.Text = "GA1-220501046-AA3-EV01"
$target.Format.Alignment = 0
$target.Format.LeftIndent = 176.4

# --- 2. Remove the now-empty paragraph that followed it ---
$afterIndex = $target.Index + 1
$empty = $d.Paragraphs.Item($afterIndex)
$empty.Range.Delete()

# --- 3. Mark the three image runs as NoProofing ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = 1
}
